$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.754.87'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '1.565.59'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.489'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.30%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.96'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").Value = '1.786.90'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '1.567.01'
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.514'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '26.808.50'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.46'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.21%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("E20").Value = '  -1.81%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  +0.62%  '
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("E31").Value = '  -3.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").Value = '1.385.50'
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("E37").Value = '  -4.03%  '
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.815'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E42").Value = '  +1.54%  '
$ws.Range("E43").Value = '  +1.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = '1.700.29'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("D49").Value = '0.0₇0986'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0951'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("E51").Value = '  -0.66%  '
